$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.500.82"
$ws.Range("E2").Value = "  +10.58%  "
$ws.Range("D3").Value = "3.252.10"
$ws.Range("E3").Value = "  +6.47%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'399.76"
$ws.Range("E5").Value = "  +2.23%  "
$ws.Range("D6").Value = "'111.07"
$ws.Range("E6").Value = "  +10.13%  "
$ws.Range("E7").Value = "  +4.89%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.619"
$ws.Range("E9").Value = "  +6.96%  "
$ws.Range("D10").Value = "'39.46"
$ws.Range("E10").Value = "  +7.84%  "
$ws.Range("D11").Value = "'0.0944"
$ws.Range("E11").Value = "  +11.71%  "
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("D13").Value = "3.768.75"
$ws.Range("E13").Value = "  +6.48%  "
$ws.Range("D14").Value = "'19.22"
$ws.Range("E14").Value = "  +5.50%  "
$ws.Range("D15").Value = "'8.10"
$ws.Range("E15").Value = "  +6.11%  "
$ws.Range("D16").Value = "3.256.82"
$ws.Range("E16").Value = "  +7.97%  "
$ws.Range("E17").Value = "  +5.22%  "
$ws.Range("D18").Value = "'10.94"
$ws.Range("E18").Value = "  +3.84%  "
$ws.Range("D19").Value = "56.453.43"
$ws.Range("E19").Value = "  +10.48%  "
$ws.Range("E20").Value = "  +5.15%  "
$ws.Range("E21").Value = "  +8.89%  "
$ws.Range("D22").Value = "'13.03"
$ws.Range("E22").Value = "  +6.63%  "
$ws.Range("D23").Value = "'299.63"
$ws.Range("E23").Value = "  +13.84%  "
$ws.Range("D24").Value = "'74.92"
$ws.Range("E24").Value = "  +7.73%  "
$ws.Range("D25").Value = "'3.22"
$ws.Range("E25").Value = "  +2.77%  "
$ws.Range("D26").Value = "'8.08"
$ws.Range("E26").Value = "  +2.92%  "
$ws.Range("D27").Value = "'28.16"
$ws.Range("E27").Value = "  +5.52%  "
$ws.Range("E28").Value = "  +5.14%  "
$ws.Range("D29").Value = "'7.33"
$ws.Range("E29").Value = "  +3.39%  "
$ws.Range("E30").Value = "  +4.58%  "
$ws.Range("D31").Value = "'0.999"
$ws.Range("D32").Value = "'0.111"
$ws.Range("E32").Value = "  +6.64%  "
$ws.Range("D33").Value = "'11.02"
$ws.Range("E33").Value = "  +5.20%  "
$ws.Range("D34").Value = "'38.52"
$ws.Range("E34").Value = "  +8.46%  "
$ws.Range("D35").Value = "'0.0488"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  +4.39%  "
$ws.Range("D37").Value = "'51.56"
$ws.Range("E37").Value = "  +3.22%  "
$ws.Range("D38").Value = "'3.13"
$ws.Range("E38").Value = "  +27.02%  "
$ws.Range("E39").Value = "  +5.25%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.12%  "
$ws.Range("D41").Value = "'17.53"
$ws.Range("E41").Value = "  +6.09%  "
$ws.Range("E42").Value = "  +6.60%  "
$ws.Range("D43").Value = "'133.68"
$ws.Range("E43").Value = "  +3.58%  "
$ws.Range("E44").Value = "  +4.80%  "
$ws.Range("D45").Value = "'3.97"
$ws.Range("E45").Value = "  +5.63%  "
$ws.Range("D46").Value = "'0.285"
$ws.Range("E46").Value = "  -1.70%  "
$ws.Range("D47").Value = "'22.12"
$ws.Range("E47").Value = "  +2.36%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.148.48"
$ws.Range("E48").Value = "  +4.18%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'2.09"
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'2.07"
$ws.Range("E50").Value = "  +46.58%  "
$ws.Range("D51").Value = "'2.42"
$ws.Range("E51").Value = "  -2.88%  "
